$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 525.4286
$ws.Range("I32").Value = 607
$ws.Range("J32").Value = 416.66666
$ws.Range("K32").Value = 607
$ws.Range("L32").Value = 416.66666
$ws.Range("M32").Value = -281
$ws.Range("N32").Value = -1068.66666

$ws.Range("H51").Value = 2000
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 2000
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -2968

$ws.Range("H98").Value = 753.46155
$ws.Range("I98").Value = 745
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 745
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = 753
$ws.Range("N98").Value = -3796

$ws.Range("H122").Value = 753.46155
$ws.Range("I122").Value = 745
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 2235
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = 215
$ws.Range("N122").Value = -7300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4269.567
$ws.Range("I74").Value = 4899.875
$ws.Range("J74").Value = 1748.3334
$ws.Range("K74").Value = 4899.875
$ws.Range("L74").Value = 1748.3334
$ws.Range("M74").Value = -4025.875
$ws.Range("N74").Value = -3496.3334

$ws.Range("H77").Value = 4269.567
$ws.Range("I77").Value = 4899.875
$ws.Range("J77").Value = 1748.3334
$ws.Range("K77").Value = 24499.375
$ws.Range("L77").Value = 8741.666999999999
$ws.Range("M77").Value = -20131.375
$ws.Range("N77").Value = -17477.667

$ws.Range("H122").Value = 1445.8889
$ws.Range("I122").Value = 1299.6666
$ws.Range("J122").Value = 1519
$ws.Range("K122").Value = 3898.9998
$ws.Range("L122").Value = 4557
$ws.Range("M122").Value = -1448.9998
$ws.Range("N122").Value = -9457

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 841.6667
$ws.Range("J64").Value = 396.66666
$ws.Range("L64").Value = 396.66666
$ws.Range("N64").Value = -846.66666

$ws.Range("H67").Value = 841.6667
$ws.Range("J67").Value = 396.66666
$ws.Range("L67").Value = 396.66666
$ws.Range("N67").Value = -1956.66666

$ws.Range("H134").Value = 3334.9412
$ws.Range("I134").Value = 3309.5
$ws.Range("K134").Value = 9928.5
$ws.Range("M134").Value = -7393.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 949.51514
$ws.Range("I31").Value = 760.55
$ws.Range("J31").Value = 1240.2307
$ws.Range("K31").Value = 760.55
$ws.Range("L31").Value = 1240.2307
$ws.Range("M31").Value = -465.55
$ws.Range("N31").Value = -1830.2307

$ws.Range("H34").Value = 949.51514
$ws.Range("I34").Value = 760.55
$ws.Range("J34").Value = 1240.2307
$ws.Range("K34").Value = 760.55
$ws.Range("L34").Value = 1240.2307
$ws.Range("M34").Value = -558.55
$ws.Range("N34").Value = -1644.2307

$ws.Range("H99").Value = 3090.1052
$ws.Range("I99").Value = 2421.2
$ws.Range("J99").Value = 3833.3333
$ws.Range("K99").Value = 2421.2
$ws.Range("L99").Value = 3833.3333
$ws.Range("M99").Value = -923.1999999999998
$ws.Range("N99").Value = -6829.3333

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 852.2857
$ws.Range("I122").Value = 864
$ws.Range("J122").Value = 843.5
$ws.Range("K122").Value = 2592
$ws.Range("L122").Value = 2530.5
$ws.Range("M122").Value = -142
$ws.Range("N122").Value = -7430.5

$ws.Range("H126").Value = 3090.1052
$ws.Range("I126").Value = 2421.2
$ws.Range("J126").Value = 3833.3333
$ws.Range("K126").Value = 7263.599999999999
$ws.Range("L126").Value = 11499.9999
$ws.Range("M126").Value = -4793.599999999999
$ws.Range("N126").Value = -16439.9999

$ws.Range("H134").Value = 3643
$ws.Range("I134").Value = 1987.7778
$ws.Range("J134").Value = 6125.8335
$ws.Range("K134").Value = 5963.3334
$ws.Range("L134").Value = 18377.5005
$ws.Range("M134").Value = -3428.3334
$ws.Range("N134").Value = -23447.5005

$ws.Range("H140").Value = 63436.863
$ws.Range("J140").Value = 63436.863
$ws.Range("L140").Value = 63436.863
$ws.Range("N140").Value = -73796.863

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 529.5
$ws.Range("I5").Value = 443.375
$ws.Range("J5").Value = 567.7778
$ws.Range("K5").Value = 1330.125
$ws.Range("L5").Value = 1703.3334
$ws.Range("M5").Value = -1218.125
$ws.Range("N5").Value = -1927.3334

$ws.Range("H109").Value = 3995.7646
$ws.Range("I109").Value = 1190.75
$ws.Range("J109").Value = 4858.846
$ws.Range("K109").Value = 3572.25
$ws.Range("L109").Value = 14576.538
$ws.Range("M109").Value = -2532.25
$ws.Range("N109").Value = -16656.538

$ws.Range("H113").Value = 1078061
$ws.Range("I113").Value = 3135228.2
$ws.Range("J113").Value = 497.2857
$ws.Range("K113").Value = 9405684.600000001
$ws.Range("L113").Value = 1491.8571
$ws.Range("M113").Value = -9403514.600000001
$ws.Range("N113").Value = -5831.8571

$ws.Range("H122").Value = 11866.263
$ws.Range("I122").Value = 30722.857
$ws.Range("J122").Value = 866.5833
$ws.Range("K122").Value = 276505.713
$ws.Range("L122").Value = 7799.2497
$ws.Range("M122").Value = -274055.713
$ws.Range("N122").Value = -12699.2497

$ws.Range("H131").Value = 4528.9644
$ws.Range("I131").Value = 753
$ws.Range("J131").Value = 4742.698
$ws.Range("K131").Value = 2259
$ws.Range("L131").Value = 14228.094
$ws.Range("M131").Value = 2781
$ws.Range("N131").Value = -24308.094

$ws.Range("H132").Value = 1123331.5
$ws.Range("J132").Value = 1444026.2
$ws.Range("L132").Value = 12996235.8
$ws.Range("N132").Value = -13001295.8

$ws.Range("H135").Value = 529.5
$ws.Range("I135").Value = 443.375
$ws.Range("J135").Value = 567.7778
$ws.Range("K135").Value = 3990.375
$ws.Range("L135").Value = 5110.000199999999
$ws.Range("M135").Value = -1455.375
$ws.Range("N135").Value = -10180.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1745.4231
$ws.Range("I102").Value = 1638.7727
$ws.Range("K102").Value = 1638.7727
$ws.Range("M102").Value = -16.77269999999999

$ws.Range("H126").Value = 2548.7
$ws.Range("I126").Value = 1689.4
$ws.Range("K126").Value = 5068.200000000001
$ws.Range("M126").Value = -2598.200000000001

$ws.Range("H132").Value = 4532.6577
$ws.Range("I132").Value = 5017.24
$ws.Range("J132").Value = 3600.7693
$ws.Range("K132").Value = 15051.72
$ws.Range("L132").Value = 10802.3079
$ws.Range("M132").Value = -12521.72
$ws.Range("N132").Value = -15862.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 294969.22
$ws.Range("I22").Value = 435596.66
$ws.Range("K22").Value = 435596.66
$ws.Range("M22").Value = -435301.66

$ws.Range("H27").Value = 294969.22
$ws.Range("I27").Value = 435596.66
$ws.Range("K27").Value = 435596.66
$ws.Range("M27").Value = -435489.66

$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2965.5557
$ws.Range("I122").Value = 2300.7778
$ws.Range("J122").Value = 3630.3333
$ws.Range("K122").Value = 6902.3334
$ws.Range("L122").Value = 10890.9999
$ws.Range("M122").Value = -4452.3334
$ws.Range("N122").Value = -15790.9999
